$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4-9 (entries beyond the two new ones that remain)
$ws.Range("A4:D9").EntireRow.Delete()

# Update row 2 with new article data
$ws.Range("B2").Value = "[PDF][PDF] Efecto del manejo forestal en la diversidad y composición arbórea de un bosque templado del noroeste de México"
$ws.Range("C2").Value = "https://www.scielo.org.mx/pdf/rcscfa/v19n2/v19n2a2.pdf"
$ws.Range("D2").Value = "En la presente investigación se evaluó el efecto de las prácticas silvícolas en la diversidad `ny composición de especies arbóreas de un bosque templado del noroeste de México. Para …"

# Update row 3 with new article data
$ws.Range("B3").Value = "[PDF][PDF] Componentes químicos y su relación con las actividades biológicas de algunos extractos vegetales"
$ws.Range("C3").Value = "https://www.redalyc.org/pdf/863/86314868005.pdf"
$ws.Range("D3").Value = "Los aceites esenciales y los extractos vegetales son mezclas complejas de metabolitos `nsecundarios que cubren un amplio espectro de efectos farmacológicos mostrando diversas …"
